$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "181703635"
$ws.Range("D2").Value = "Imported93496618"
$ws.Range("P2").Value = "Sridevi73@zenwork.com"
